$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 'last updated' timestamp in A1
$ws.Range("A1").Value = 'Datos actualizados a 8 de Octubre de 2020 a las 17:53'

# Update country rows whose data/ranking changed
$ws.Cells.Item(4, 1).Value = 'Estados Unidos'
$ws.Cells.Item(4, 2).Value = 7786429
$ws.Cells.Item(4, 3).Value = 10205
$ws.Cells.Item(4, 4).Value = 4997380
$ws.Cells.Item(4, 5).Value = 2571983
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 282
$ws.Cells.Item(4, 8).Value = 217066

$ws.Cells.Item(5, 1).Value = 'India'
$ws.Cells.Item(5, 2).Value = 6841813
$ws.Cells.Item(5, 3).Value = 8825
$ws.Cells.Item(5, 4).Value = 5836826
$ws.Cells.Item(5, 5).Value = 899355
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 78
$ws.Cells.Item(5, 8).Value = 105632

$ws.Cells.Item(15, 1).Value = 'Reino Unido'
$ws.Cells.Item(15, 2).Value = 561815
$ws.Cells.Item(15, 3).Value = 17540
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 77
$ws.Cells.Item(15, 8).Value = 42592

$ws.Cells.Item(17, 1).Value = 'Chile'
$ws.Cells.Item(17, 2).Value = 476016
$ws.Cells.Item(17, 3).Value = 1576
$ws.Cells.Item(17, 4).Value = 448710
$ws.Cells.Item(17, 5).Value = 14139
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 77
$ws.Cells.Item(17, 8).Value = 13167

$ws.Cells.Item(20, 1).Value = 'Italia'
$ws.Cells.Item(20, 2).Value = 338398
$ws.Cells.Item(20, 3).Value = 4458
$ws.Cells.Item(20, 4).Value = 236363
$ws.Cells.Item(20, 5).Value = 65952
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 22
$ws.Cells.Item(20, 8).Value = 36083

$ws.Cells.Item(21, 1).Value = 'Arabia Saudita'
$ws.Cells.Item(21, 2).Value = 338132
$ws.Cells.Item(21, 3).Value = 421
$ws.Cells.Item(21, 4).Value = 323769
$ws.Cells.Item(21, 5).Value = 9391
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 25
$ws.Cells.Item(21, 8).Value = 4972

$ws.Cells.Item(26, 1).Value = 'Alemania'
$ws.Cells.Item(26, 2).Value = 312679
$ws.Cells.Item(26, 3).Value = 1566
$ws.Cells.Item(26, 4).Value = 267700
$ws.Cells.Item(26, 5).Value = 35315
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).Value = 12
$ws.Cells.Item(26, 8).Value = 9664

$ws.Cells.Item(29, 1).Value = 'Canada'
$ws.Cells.Item(29, 2).Value = 173920
$ws.Cells.Item(29, 3).Value = 797
$ws.Cells.Item(29, 4).Value = 146361
$ws.Cells.Item(29, 5).Value = 18014
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(29, 7).Value = 4
$ws.Cells.Item(29, 8).Value = 9545

$ws.Cells.Item(38, 1).Value = 'Republica Dominicana'
$ws.Cells.Item(38, 2).Value = 116872
$ws.Cells.Item(38, 3).Value = 724
$ws.Cells.Item(38, 4).Value = 92567
$ws.Cells.Item(38, 5).Value = 22142
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 4
$ws.Cells.Item(38, 8).Value = 2163

$ws.Cells.Item(47, 1).Value = 'Guatemala'
$ws.Cells.Item(47, 2).Value = 96480
$ws.Cells.Item(47, 3).Value = 776
$ws.Cells.Item(47, 4).Value = 84738
$ws.Cells.Item(47, 5).Value = 8395
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(47, 7).Value = 12
$ws.Cells.Item(47, 8).Value = 3347

$ws.Cells.Item(49, 1).Value = 'Japon'
$ws.Cells.Item(49, 2).Value = 87020
$ws.Cells.Item(49, 3).Value = 477
$ws.Cells.Item(49, 4).Value = 80227
$ws.Cells.Item(49, 5).Value = 5180
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(49, 7).Value = 8
$ws.Cells.Item(49, 8).Value = 1613

$ws.Cells.Item(58, 1).Value = 'Moldavia'
$ws.Cells.Item(58, 2).Value = 59915
$ws.Cells.Item(58, 3).Value = 1121
$ws.Cells.Item(58, 4).Value = 43008
$ws.Cells.Item(58, 5).Value = 15483
$ws.Cells.Item(58, 6).Value = 0
$ws.Cells.Item(58, 7).Value = 18
$ws.Cells.Item(58, 8).Value = 1424

$ws.Cells.Item(59, 1).Value = 'Uzbekistan'
$ws.Cells.Item(59, 2).Value = 59905
$ws.Cells.Item(59, 3).Value = 326
$ws.Cells.Item(59, 4).Value = 56568
$ws.Cells.Item(59, 5).Value = 2843
$ws.Cells.Item(59, 6).Value = 0
$ws.Cells.Item(59, 7).Value = 3
$ws.Cells.Item(59, 8).Value = 494

$ws.Cells.Item(60, 1).Value = 'Nigeria'
$ws.Cells.Item(60, 2).Value = 59738
$ws.Cells.Item(60, 3).Value = 0
$ws.Cells.Item(60, 4).Value = 51403
$ws.Cells.Item(60, 5).Value = 7222
$ws.Cells.Item(60, 6).Value = 0
$ws.Cells.Item(60, 7).Value = 0
$ws.Cells.Item(60, 8).Value = 1113

$ws.Cells.Item(61, 1).Value = 'Suiza'
$ws.Cells.Item(61, 2).Value = 58881
$ws.Cells.Item(61, 3).Value = 1172
$ws.Cells.Item(61, 4).Value = 47300
$ws.Cells.Item(61, 5).Value = 9494
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 5
$ws.Cells.Item(61, 8).Value = 2087

$ws.Cells.Item(62, 1).Value = 'Singapur'
$ws.Cells.Item(62, 2).Value = 57849
$ws.Cells.Item(62, 3).Value = 9
$ws.Cells.Item(62, 4).Value = 57668
$ws.Cells.Item(62, 5).Value = 154
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(62, 8).Value = 27

$ws.Cells.Item(86, 1).Value = 'Grecia'
$ws.Cells.Item(86, 2).Value = 21381
$ws.Cells.Item(86, 3).Value = 434
$ws.Cells.Item(86, 4).Value = 9989
$ws.Cells.Item(86, 5).Value = 10962
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 6
$ws.Cells.Item(86, 8).Value = 430

$ws.Cells.Item(100, 1).Value = 'Montenegro'
$ws.Cells.Item(100, 2).Value = 13004
$ws.Cells.Item(100, 3).Value = 210
$ws.Cells.Item(100, 4).Value = 9154
$ws.Cells.Item(100, 5).Value = 3659
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 7).Value = 1
$ws.Cells.Item(100, 8).Value = 191

$ws.Cells.Item(107, 1).Value = 'Guayana Francesa'
$ws.Cells.Item(107, 2).Value = 10128
$ws.Cells.Item(107, 3).Value = 25
$ws.Cells.Item(107, 4).Value = 9799
$ws.Cells.Item(107, 5).Value = 260
$ws.Cells.Item(107, 6).Value = 0
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = 69

$ws.Cells.Item(117, 1).Value = 'Jamaica'
$ws.Cells.Item(117, 2).Value = 7273
$ws.Cells.Item(117, 3).Value = 82
$ws.Cells.Item(117, 4).Value = 2732
$ws.Cells.Item(117, 5).Value = 4413
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 7).Value = 2
$ws.Cells.Item(117, 8).Value = 128

$ws.Cells.Item(120, 1).Value = 'Cuba'
$ws.Cells.Item(120, 2).Value = 5917
$ws.Cells.Item(120, 3).Value = 19
$ws.Cells.Item(120, 4).Value = 5371
$ws.Cells.Item(120, 5).Value = 423
$ws.Cells.Item(120, 6).Value = 0
$ws.Cells.Item(120, 7).Value = 0
$ws.Cells.Item(120, 8).Value = 123

$ws.Cells.Item(136, 1).Value = 'Sri Lanka'
$ws.Cells.Item(136, 2).Value = 4469
$ws.Cells.Item(136, 3).Value = 10
$ws.Cells.Item(136, 4).Value = 3278
$ws.Cells.Item(136, 5).Value = 1178
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 13

$ws.Cells.Item(150, 1).Value = 'Principado de Andorra'
$ws.Cells.Item(150, 2).Value = 2568
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(150, 4).Value = 1715
$ws.Cells.Item(150, 5).Value = 799
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 1
$ws.Cells.Item(150, 8).Value = 54

$ws.Cells.Item(151, 1).Value = 'Polinesia Francesa'
$ws.Cells.Item(151, 2).Value = 2420
$ws.Cells.Item(151, 3).Value = 62
$ws.Cells.Item(151, 4).Value = 1857
$ws.Cells.Item(151, 5).Value = 553
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 10

$ws.Cells.Item(152, 1).Value = 'Benin'
$ws.Cells.Item(152, 2).Value = 2411
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(152, 4).Value = 1973
$ws.Cells.Item(152, 5).Value = 397
$ws.Cells.Item(152, 6).Value = 0
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 41

$ws.Cells.Item(153, 1).Value = 'Guinea-Bisau'
$ws.Cells.Item(153, 2).Value = 2385
$ws.Cells.Item(153, 3).Value = 0
$ws.Cells.Item(153, 4).Value = 1728
$ws.Cells.Item(153, 5).Value = 617
$ws.Cells.Item(153, 6).Value = 0
$ws.Cells.Item(153, 7).Value = 0
$ws.Cells.Item(153, 8).Value = 40

$ws.Cells.Item(154, 1).Value = 'Letonia'
$ws.Cells.Item(154, 2).Value = 2370
$ws.Cells.Item(154, 3).Value = 109
$ws.Cells.Item(154, 4).Value = 1322
$ws.Cells.Item(154, 5).Value = 1008
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 40

$ws.Cells.Item(160, 1).Value = 'Republica de Chipre'
$ws.Cells.Item(160, 2).Value = 1918
$ws.Cells.Item(160, 3).Value = 21
$ws.Cells.Item(160, 4).Value = 1369
$ws.Cells.Item(160, 5).Value = 525
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 0
$ws.Cells.Item(160, 8).Value = 24

$ws.Cells.Item(161, 1).Value = 'Togo'
$ws.Cells.Item(161, 2).Value = 1898
$ws.Cells.Item(161, 3).Value = 0
$ws.Cells.Item(161, 4).Value = 1419
$ws.Cells.Item(161, 5).Value = 430
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = 0
$ws.Cells.Item(161, 8).Value = 49

$ws.Cells.Item(164, 1).Value = 'Lesoto'
$ws.Cells.Item(164, 2).Value = 1786
$ws.Cells.Item(164, 3).Value = 19
$ws.Cells.Item(164, 4).Value = 926
$ws.Cells.Item(164, 5).Value = 819
$ws.Cells.Item(164, 6).Value = 0
$ws.Cells.Item(164, 7).Value = 1
$ws.Cells.Item(164, 8).Value = 41

$ws.Cells.Item(170, 1).Value = 'San Marino'
$ws.Cells.Item(170, 2).Value = 741
$ws.Cells.Item(170, 3).Value = 9
$ws.Cells.Item(170, 4).Value = 682
$ws.Cells.Item(170, 5).Value = 17
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 7).Value = 0
$ws.Cells.Item(170, 8).Value = 42

$ws.Cells.Item(174, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(174, 2).Value = 549
$ws.Cells.Item(174, 3).Value = 8
$ws.Cells.Item(174, 4).Value = 532
$ws.Cells.Item(174, 5).Value = 10
$ws.Cells.Item(174, 6).Value = 0
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 7

$ws.Cells.Item(183, 1).Value = 'Eritrea'
$ws.Cells.Item(183, 2).Value = 405
$ws.Cells.Item(183, 3).Value = 7
$ws.Cells.Item(183, 4).Value = 364
$ws.Cells.Item(183, 5).Value = 41
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 0

$ws.Cells.Item(215, 1).Value = 'Montserrat'
$ws.Cells.Item(215, 2).Value = 13
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 12
$ws.Cells.Item(215, 5).Value = 0
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 1

$ws.Cells.Item(216, 1).Value = 'Islas Malvinas'
$ws.Cells.Item(216, 2).Value = 13
$ws.Cells.Item(216, 3).Value = 0
$ws.Cells.Item(216, 4).Value = 13
$ws.Cells.Item(216, 5).Value = 0
$ws.Cells.Item(216, 6).Value = 0
$ws.Cells.Item(216, 7).Value = 0
$ws.Cells.Item(216, 8).Value = 0
